$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that needs to be bumped
# by one day (45184 -> 45185) for every data row (rows 2 through 13).
$range = $ws.Range("C2:C13")
foreach ($cell in $range.Cells) {
    $cell.Value = $cell.Value2 + 1
}
